$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "depth" (H) and "depth_std" (I) columns for rows 2-8
$ws.Range("H2").Value = 0.5297855048019494
$ws.Range("I2").Value = 0.05306360429042664

$ws.Range("H3").Value = 0.9975627310283326
$ws.Range("I3").Value = 0.1155566321558266

$ws.Range("H4").Value = 1.065315069088294
$ws.Range("I4").Value = 0.1427956937476589

$ws.Range("H5").Value = 0.7470287814537109
$ws.Range("I5").Value = 0.1677466073687303

$ws.Range("H6").Value = 0.7171265693569835
$ws.Range("I6").Value = 0.2515030733153356

$ws.Range("H7").Value = 0.8166802768446795
$ws.Range("I7").Value = 0.257224835533526

$ws.Range("H8").Value = 0.8688388354752572
$ws.Range("I8").Value = 0.2954030563117702
